# Discussion points 08 Apr 2020.docx - edit script
#
# 1. Insert four new bullet paragraphs (Service booking Ref No. / While
#    selecting ref no / Date to be auto populated / Actual booking date)
#    right after the "Driver Name" bullet, under the "Operations -
#    Workshop Admission" heading.
# 2. Move <w:lastRenderedPageBreak/> from the "When CLOSING a job card..."
#    run to the "Give a textbox to enter the Inspection Ref No..." run.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: insert the four new paragraphs after the "Driver Name" bullet
# ---------------------------------------------------------------------
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "ehicle to workshop and that list will not be there in the system.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'Driver Name' bullet paragraph."
}

$insertPos = $anchor.End
$insertRange = $d.Range($insertPos, $insertPos)

$newParagraphsXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Service booking Ref No. </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&#8211;</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Show all pending bookings</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>While selecting ref no, show vehicle details</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Date to be auto populated with current date</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Actual booking date &amp; Actual admission date to be captured.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($newParagraphsXml)

# ---------------------------------------------------------------------
# Step 2a: add <w:lastRenderedPageBreak/> to the "Give a textbox..." run
# ---------------------------------------------------------------------
$giveRange = $d.Content
$found = $giveRange.Find.Execute(
    "Give a textbox to enter the Inspection Ref No",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'Give a textbox...' paragraph."
}

$giveRange.Expand(4)  # wdParagraph

$giveParagraphXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:ind w:left="1440"/>
    <w:jc w:val="both"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Give a textbox to enter the Inspection Ref No. &amp; another textbox to enter the observations.</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve"> A File Upload control to upload the scanned copy of Inspection sheet (Actually Inspection will be done of paper based forms). This attachment should be allowed even after closing the job card.</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$giveRange.InsertXML($giveParagraphXml)

# ---------------------------------------------------------------------
# Step 2b: remove <w:lastRenderedPageBreak/> from the "When CLOSING..." run
# ---------------------------------------------------------------------
$closingRange = $d.Content
$found = $closingRange.Find.Execute(
    "When CLOSING a job card, there should be a popup",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'When CLOSING a job card...' paragraph."
}

$closingRange.Expand(4)  # wdParagraph

$closingParagraphXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:spacing w:line="240" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve">When CLOSING a job card, there should be a popup that needs to be filled by the technician. </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
    <w:t>We will decide w</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
    <w:t xml:space="preserve">hat is to be there in </w:t>
  </w:r>
  <w:r>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
    <w:t>the popup will be decided later.</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$closingRange.InsertXML($closingParagraphXml)

Write-Output "Edit complete."
